# Adds the 33 new MIMS.me.air v4.0 columns (N:AT) to the header row (row 15),
# each with its harmonized-name value, a descriptive cell comment, and the
# same 'optional/yellow' header style already used for the existing columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C15 already carries the yellow/optional header style (s="8") that the new
# header cells should use; copy its format onto each new header cell below.
$styleSource = $ws.Range("C15")

$ws.Range("N15").Value = 'barometric_press'
$ws.Range("N15").AddComment('force per unit area exerted against a surface by the weight of air above that surface')
$styleSource.Copy()
$ws.Range("N15").PasteSpecial(-4122)

$ws.Range("O15").Value = 'carb_dioxide'
$ws.Range("O15").AddComment('carbon dioxide (gas) amount or concentration at the time of sampling')
$styleSource.Copy()
$ws.Range("O15").PasteSpecial(-4122)

$ws.Range("P15").Value = 'carb_monoxide'
$ws.Range("P15").AddComment('carbon monoxide (gas) amount or concentration at the time of sampling')
$styleSource.Copy()
$ws.Range("P15").PasteSpecial(-4122)

$ws.Range("Q15").Value = 'chem_administration'
$ws.Range("Q15").AddComment('list of chemical compounds administered to the host or site where sampling occurred, and when (e.g. antibiotics, N fertilizer, air filter); can include multiple compounds. For Chemical Entities of Biological Interest ontology (CHEBI) (v1.72), please see http://bioportal.bioontology.org/visualize/44603')
$styleSource.Copy()
$ws.Range("Q15").PasteSpecial(-4122)

$ws.Range("R15").Value = 'elev'
$ws.Range("R15").AddComment('The elevation of the sampling site as measured by the vertical distance from mean sea level.')
$styleSource.Copy()
$ws.Range("R15").PasteSpecial(-4122)

$ws.Range("S15").Value = 'humidity'
$ws.Range("S15").AddComment('amount of water vapour in the air, at the time of sampling')
$styleSource.Copy()
$ws.Range("S15").PasteSpecial(-4122)

$ws.Range("T15").Value = 'isolation_source'
$ws.Range("T15").AddComment('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.')
$styleSource.Copy()
$ws.Range("T15").PasteSpecial(-4122)

$ws.Range("U15").Value = 'methane'
$ws.Range("U15").AddComment('methane (gas) amount or concentration at the time of sampling')
$styleSource.Copy()
$ws.Range("U15").PasteSpecial(-4122)

$ws.Range("V15").Value = 'misc_param'
$ws.Range("V15").AddComment('any other measurement performed or parameter collected, that is not listed here')
$styleSource.Copy()
$ws.Range("V15").PasteSpecial(-4122)

$ws.Range("W15").Value = 'organism_count'
$ws.Range("W15").AddComment('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts')
$styleSource.Copy()
$ws.Range("W15").PasteSpecial(-4122)

$ws.Range("X15").Value = 'oxy_stat_samp'
$ws.Range("X15").AddComment('oxygenation status of sample')
$styleSource.Copy()
$ws.Range("X15").PasteSpecial(-4122)

$ws.Range("Y15").Value = 'oxygen'
$ws.Range("Y15").AddComment('oxygen (gas) amount or concentration at the time of sampling')
$styleSource.Copy()
$ws.Range("Y15").PasteSpecial(-4122)

$ws.Range("Z15").Value = 'perturbation'
$ws.Range("Z15").AddComment('type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types')
$styleSource.Copy()
$ws.Range("Z15").PasteSpecial(-4122)

$ws.Range("AA15").Value = 'pollutants'
$ws.Range("AA15").AddComment('pollutant types and, amount or concentrations measured at the time of sampling; can report multiple pollutants by entering numeric values preceded by name of pollutant')
$styleSource.Copy()
$ws.Range("AA15").PasteSpecial(-4122)

$ws.Range("AB15").Value = 'ref_biomaterial'
$ws.Range("AB15").AddComment('Primary publication or genome report in the form of pubmed ID, DOI or URL')
$styleSource.Copy()
$ws.Range("AB15").PasteSpecial(-4122)

$ws.Range("AC15").Value = 'rel_to_oxygen'
$ws.Range("AC15").AddComment('Aerobic or anaerobic')
$styleSource.Copy()
$ws.Range("AC15").PasteSpecial(-4122)

$ws.Range("AD15").Value = 'resp_part_matter'
$ws.Range("AD15").AddComment('concentration of substances that remain suspended in the air, and comprise mixtures of organic and inorganic substances (PM10 and PM2.5); can report multiple PM''s by entering numeric values preceded by name of PM')
$styleSource.Copy()
$ws.Range("AD15").PasteSpecial(-4122)

$ws.Range("AE15").Value = 'samp_collect_device'
$ws.Range("AE15").AddComment('Method or device employed for collecting sample')
$styleSource.Copy()
$ws.Range("AE15").PasteSpecial(-4122)

$ws.Range("AF15").Value = 'samp_mat_process'
$ws.Range("AF15").AddComment('Processing applied to the sample during or after isolation')
$styleSource.Copy()
$ws.Range("AF15").PasteSpecial(-4122)

$ws.Range("AG15").Value = 'samp_salinity'
$ws.Range("AG15").AddComment('salinity of sample, i.e. measure of total salt concentration')
$styleSource.Copy()
$ws.Range("AG15").PasteSpecial(-4122)

$ws.Range("AH15").Value = 'samp_size'
$ws.Range("AH15").AddComment('Amount or size of sample (volume, mass or area) that was collected')
$styleSource.Copy()
$ws.Range("AH15").PasteSpecial(-4122)

$ws.Range("AI15").Value = 'samp_store_dur'
$ws.Range("AI15").AddComment('duration for which sample was stored')
$styleSource.Copy()
$ws.Range("AI15").PasteSpecial(-4122)

$ws.Range("AJ15").Value = 'samp_store_loc'
$ws.Range("AJ15").AddComment('location at which sample was stored, usually name of a specific freezer/room')
$styleSource.Copy()
$ws.Range("AJ15").PasteSpecial(-4122)

$ws.Range("AK15").Value = 'samp_store_temp'
$ws.Range("AK15").AddComment('temperature at which sample was stored, e.g. -80')
$styleSource.Copy()
$ws.Range("AK15").PasteSpecial(-4122)

$ws.Range("AL15").Value = 'samp_vol_we_dna_ext'
$ws.Range("AL15").AddComment('volume (mL) or weight (g) of sample processed for DNA extraction')
$styleSource.Copy()
$ws.Range("AL15").PasteSpecial(-4122)

$ws.Range("AM15").Value = 'solar_irradiance'
$ws.Range("AM15").AddComment('the amount of solar energy that arrives at a specific area of a surface during a specific time interval')
$styleSource.Copy()
$ws.Range("AM15").PasteSpecial(-4122)

$ws.Range("AN15").Value = 'source_material_id'
$ws.Range("AN15").AddComment('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.')
$styleSource.Copy()
$ws.Range("AN15").PasteSpecial(-4122)

$ws.Range("AO15").Value = 'temperature'
$ws.Range("AO15").AddComment('temperature of the sample at time of sampling')
$styleSource.Copy()
$ws.Range("AO15").PasteSpecial(-4122)

$ws.Range("AP15").Value = 'ventilation_rate'
$ws.Range("AP15").AddComment('ventilation rate of the system in the sampled premises')
$styleSource.Copy()
$ws.Range("AP15").PasteSpecial(-4122)

$ws.Range("AQ15").Value = 'ventilation_type'
$ws.Range("AQ15").AddComment('ventilation system used in the sampled premises')
$styleSource.Copy()
$ws.Range("AQ15").PasteSpecial(-4122)

$ws.Range("AR15").Value = 'volatile_org_comp'
$ws.Range("AR15").AddComment('concentration of carbon-based chemicals that easily evaporate at room temperature; can report multiple volatile organic compounds by entering numeric values preceded by name of compound')
$styleSource.Copy()
$ws.Range("AR15").PasteSpecial(-4122)

$ws.Range("AS15").Value = 'wind_direction'
$ws.Range("AS15").AddComment('wind direction is the direction from which a wind originates')
$styleSource.Copy()
$ws.Range("AS15").PasteSpecial(-4122)

$ws.Range("AT15").Value = 'wind_speed'
$ws.Range("AT15").AddComment('speed of wind measured at the time of sampling')
$styleSource.Copy()
$ws.Range("AT15").PasteSpecial(-4122)

$excel.CutCopyMode = $false
